# Update Pinjaman dan Simpanan
# Inserts two new columns of data (E = "Jenis Pinjaman", H = "Bulan (Angka)")
# into the existing Saldo Pinjaman table, shifting the old E/F/G/H columns
# right into F/G/I/J, and updates the "CICILAN BULAN ..." label text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 1 (header row) - move existing headers right, two at a time, starting
# from the rightmost column so we never overwrite a cell before reading it.
# ---------------------------------------------------------------------------

# G1 "Tanggal Bayar" (s1) -> I1
$ws.Range("G1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "Tanggal Bayar"
$ws.Range("G1").ClearContents()

# H1 (empty, s3) -> J1
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("H1").ClearContents()

# F1 "Sisa Cicilan" (s1) -> G1
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Sisa Cicilan"
$ws.Range("F1").ClearContents()

# E1 "Cicilan" (s1) -> F1
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Cicilan"
$ws.Range("E1").ClearContents()

# New headers: E1 "Jenis Pinjaman", H1 "Bulan (Angka)" (style like other s1 header cells)
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Jenis Pinjaman"

$ws.Range("D1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Bulan (Angka)"

# ---------------------------------------------------------------------------
# Row 2 (data row)
# ---------------------------------------------------------------------------

# H2 rich-text note -> J2 (value only; move via cut so the rich run survives)
$ws.Range("H2").Cut()
$ws.Range("J2").PasteSpecial(-4163)

# G2 "=TODAY()" formula (s9) -> I2
$ws.Range("G2").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("I2").Formula = "=TODAY()"
$ws.Range("G2").ClearContents()

# F2 9000000 (s5) -> G2
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G2").Value = 9000000
$ws.Range("F2").ClearContents()

# E2 1000000 (s5) -> F2
$ws.Range("E2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F2").Value = 1000000
$ws.Range("E2").ClearContents()

# New data: E2 "PINJAMAN UANG", H2 9 (Bulan Angka = September)
$ws.Range("C2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("E2").Value = "PINJAMAN UANG"

$ws.Range("C2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("H2").Value = 9

# Update the "Keterangan" label for the new month
$ws.Range("D2").Value = "CICILAN BULAN SEPTEMBER 2025"

# ---------------------------------------------------------------------------
# Row 3 (border/spacer row)
# ---------------------------------------------------------------------------

# H3 "Mulai Dari Sini" -> J3
$ws.Range("H3").Cut()
$ws.Range("J3").PasteSpecial(-4163)

# G3 (empty, s7) -> I3
$ws.Range("G3").Copy()
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("G3").ClearContents()

# New blank cells G3, H3 matching the plain row-3 border style (copy from F3)
$ws.Range("F3").Copy()
$ws.Range("G3").PasteSpecial(-4122)

$ws.Range("F3").Copy()
$ws.Range("H3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Column widths (D grew to fit the longer label, E was hand-resized for the
# new "Jenis Pinjaman" column, H/I are the newly typed columns)
# ---------------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 28.166666666666668
$ws.Columns("E").ColumnWidth = 27.276041666666668
$ws.Columns("H").ColumnWidth = 11.166666666666666
$ws.Columns("I").ColumnWidth = 11.498697916666666

# ---------------------------------------------------------------------------
# Housekeeping to mirror the rest of the diff
# ---------------------------------------------------------------------------
$ws.Range("E6").Select()
